# Update "想去人数" (interest count) figures in the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 6877
$wsExhibit.Range("F4").Value  = 446
$wsExhibit.Range("F5").Value  = 74
$wsExhibit.Range("F6").Value  = 14
$wsExhibit.Range("F8").Value  = 113
$wsExhibit.Range("F11").Value = 7
$wsExhibit.Range("F12").Value = 44
$wsExhibit.Range("F15").Value = 19
$wsExhibit.Range("F16").Value = 1798
$wsExhibit.Range("F17").Value = 33
$wsExhibit.Range("F18").Value = 3489
$wsExhibit.Range("F21").Value = 21
$wsExhibit.Range("F22").Value = 2119
$wsExhibit.Range("F23").Value = 202
$wsExhibit.Range("F25").Value = 31
$wsExhibit.Range("F28").Value = 13

# Sheet "全部类型" (sheetId 4) - same underlying rows, shifted by the extra
# "演出" row inserted at row 7, so row numbers differ from "展览" sheet.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6877
$wsAll.Range("F4").Value  = 446
$wsAll.Range("F5").Value  = 74
$wsAll.Range("F6").Value  = 14
$wsAll.Range("F9").Value  = 113
$wsAll.Range("F12").Value = 7
$wsAll.Range("F13").Value = 44
$wsAll.Range("F16").Value = 19
$wsAll.Range("F17").Value = 1798
$wsAll.Range("F18").Value = 33
$wsAll.Range("F19").Value = 3489
$wsAll.Range("F22").Value = 21
$wsAll.Range("F23").Value = 2119
$wsAll.Range("F24").Value = 202
$wsAll.Range("F26").Value = 31
$wsAll.Range("F29").Value = 13
